$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B (the earliest-year data column for each species block).
# This shifts all subsequent columns (C..K) one position to the left,
# reproducing the "drop earliest year, shift table" edit described by the diff:
#   BEAR:      2014..2023 -> 2015..2023
#   WOLVERINE: 2015..2024 -> 2016..2024
#   WOLF:      2014/2015..2023/2024 -> 2015/2016..2023/2024
$ws.Columns("B").Delete()

$ws.Range("D17").Select() | Out-Null
